$d = $word.ActiveDocument

function Replace-InParagraph($index, $oldText, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# 1. Title heading (Heading1, paragraph 1)
Replace-InParagraph 1 "Play Genie's Palace for Free - High Potential Rewards" "Play Genie’s Palace for Free"

# 2. "What we like" bullet list (paragraphs 42, 43, 45)
Replace-InParagraph 42 "High potential rewards" "Traditional slot machine structure"
Replace-InParagraph 43 "Free demo version available" "High volatility with rewarding wins"
Replace-InParagraph 45 "Traditional slot machine structure" "Demo version available for free play"

# 3. "What we don't like" bullet list (paragraphs 47, 48)
Replace-InParagraph 47 "Very high volatility" "Potential for significant number of losses"
Replace-InParagraph 48 "Limited paylines compared to other slots" "Limited paylines"

# 4. Bold meta title (paragraph 49)
Replace-InParagraph 49 "Play Genie's Palace for Free - High Potential Rewards" "Play Genie’s Palace for Free"

# 5. Italic meta description (paragraph 50)
Replace-InParagraph 50 "Read our review of Genie's Palace online slot game. Play for free with high potential rewards and special bonuses. Try demo version without registration." "Explore the features and gameplay of Genie’s Palace in this review. Play for free without registration!"

Write-Output "Done"
